$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 117, pushing the existing rows
# 117-147 down to 119-149 (same as the weekly report picking up a new
# date at the top of this produce's price history).
$ws.Range("A117:A118").EntireRow.Insert()

$fecha = Get-Date -Year 2023 -Month 11 -Day 14 -Hour 0 -Minute 0 -Second 0

# Row 117: new "Primera" quality record for 2023-11-14
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = $fecha
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = 100112044
$ws.Range("G117").Value = "Perejil"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 300
$ws.Range("K117").Value = 2000
$ws.Range("L117").Value = 2000
$ws.Range("M117").Value = 2000
$ws.Range("N117").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O117").Value = "Región de Ñuble"
$ws.Range("P117").Value = 2000
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"

# Row 118: new "Segunda" quality record for 2023-11-14
$ws.Range("A118").Value = 7
$ws.Range("B118").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C118").Value = "Ñuble"
$ws.Range("D118").Value = $fecha
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = 100112044
$ws.Range("G118").Value = "Perejil"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Segunda"
$ws.Range("J118").Value = 300
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 1500
$ws.Range("M118").Value = 1500
$ws.Range("N118").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O118").Value = "Región de Ñuble"
$ws.Range("P118").Value = 1500
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = "Hortaliza"

Write-Host "Inserted rows 117-118; new dimension:" $ws.UsedRange.Address()
